$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2233
$ws.Range("J6").Value = 2599.5
$ws.Range("L6").Value = 7798.5
$ws.Range("N6").Value = -8022.5
$ws.Range("H17").Value = 2321.75
$ws.Range("J17").Value = 2321.75
$ws.Range("L17").Value = 6965.25
$ws.Range("N17").Value = -7301.25
$ws.Range("H33").Value = 483.08334
$ws.Range("I33").Value = 453.92307
$ws.Range("J33").Value = 517.5454999999999
$ws.Range("K33").Value = 453.92307
$ws.Range("L33").Value = 517.5454999999999
$ws.Range("M33").Value = -224.92307
$ws.Range("N33").Value = -975.5454999999999
$ws.Range("H51").Value = 17864326
$ws.Range("I51").Value = 62501996
$ws.Range("J51").Value = 9257.799999999999
$ws.Range("K51").Value = 62501996
$ws.Range("L51").Value = 9257.799999999999
$ws.Range("M51").Value = -62501512
$ws.Range("N51").Value = -10225.8
$ws.Range("H107").Value = 1438.125
$ws.Range("I107").Value = 1266.1111
$ws.Range("J107").Value = 1954.1666
$ws.Range("K107").Value = 1266.1111
$ws.Range("L107").Value = 1954.1666
$ws.Range("M107").Value = 653.8888999999999
$ws.Range("N107").Value = -5794.1666
$ws.Range("H116").Value = 15205.333
$ws.Range("J116").Value = 9034.416999999999
$ws.Range("L116").Value = 9034.416999999999
$ws.Range("N116").Value = -15918.417
$ws.Range("H131").Value = 4528.0835
$ws.Range("I131").Value = 4533.9
$ws.Range("J131").Value = 4499
$ws.Range("K131").Value = 13601.7
$ws.Range("L131").Value = 13497
$ws.Range("M131").Value = -8561.699999999999
$ws.Range("N131").Value = -23577
$ws.Range("H132").Value = 2761926.5
$ws.Range("I132").Value = 3055601.8
$ws.Range("K132").Value = 9166805.399999999
$ws.Range("M132").Value = -9164275.399999999
$ws.Range("H137").Value = 21081.176
$ws.Range("I137").Value = 32759.9
$ws.Range("J137").Value = 4397.2856
$ws.Range("K137").Value = 98279.70000000001
$ws.Range("L137").Value = 13191.8568
$ws.Range("M137").Value = -95729.70000000001
$ws.Range("N137").Value = -18291.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3053.4167
$ws.Range("I45").Value = 1798.7646
$ws.Range("K45").Value = 1798.7646
$ws.Range("M45").Value = -1421.7646
$ws.Range("H61").Value = 8826.066000000001
$ws.Range("I61").Value = 1479.4
$ws.Range("K61").Value = 1479.4
$ws.Range("M61").Value = -1267.4
$ws.Range("H74").Value = 278274.38
$ws.Range("I74").Value = 334307.94
$ws.Range("K74").Value = 334307.94
$ws.Range("M74").Value = -333433.94
$ws.Range("H77").Value = 278274.38
$ws.Range("I77").Value = 334307.94
$ws.Range("K77").Value = 1671539.7
$ws.Range("M77").Value = -1667171.7
$ws.Range("H122").Value = 1715.6
$ws.Range("I122").Value = 1504
$ws.Range("K122").Value = 4512
$ws.Range("M122").Value = -2062
$ws.Range("H136").Value = 8826.066000000001
$ws.Range("I136").Value = 1479.4
$ws.Range("K136").Value = 4438.200000000001
$ws.Range("M136").Value = -1888.200000000001
$ws.Range("H140").Value = 91000
$ws.Range("J140").Value = 91000
$ws.Range("L140").Value = 91000
$ws.Range("N140").Value = -101360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 52512.1
$ws.Range("I107").Value = 71901.14
$ws.Range("K107").Value = 71901.14
$ws.Range("M107").Value = -69981.14
$ws.Range("H134").Value = 2187.6829
$ws.Range("I134").Value = 1765.0303
$ws.Range("J134").Value = 3931.125
$ws.Range("K134").Value = 5295.090899999999
$ws.Range("L134").Value = 11793.375
$ws.Range("M134").Value = -2760.090899999999
$ws.Range("N134").Value = -16863.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5886607.5
$ws.Range("J31").Value = 4728.2666
$ws.Range("L31").Value = 4728.2666
$ws.Range("N31").Value = -5318.2666
$ws.Range("H34").Value = 5886607.5
$ws.Range("J34").Value = 4728.2666
$ws.Range("L34").Value = 4728.2666
$ws.Range("N34").Value = -5132.2666
$ws.Range("H50").Value = 36624.5
$ws.Range("J50").Value = 36624.5
$ws.Range("L50").Value = 36624.5
$ws.Range("N50").Value = -37874.5
$ws.Range("H100").Value = 86191.28999999999
$ws.Range("J100").Value = 86191.28999999999
$ws.Range("L100").Value = 86191.28999999999
$ws.Range("N100").Value = -88355.28999999999
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""
$ws.Range("H134").Value = 2091.7368
$ws.Range("I134").Value = 1422.6154
$ws.Range("K134").Value = 4267.8462
$ws.Range("M134").Value = -1732.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 436.66666
$ws.Range("H26").Value = 406.44446
$ws.Range("I26").Value = 344.14285
$ws.Range("K26").Value = 1032.42855
$ws.Range("M26").Value = -744.4285500000001
$ws.Range("H39").Value = 13336
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = ""
$ws.Range("H56").Value = 6874.5
$ws.Range("I56").Value = 6874.5
$ws.Range("K56").Value = 6874.5
$ws.Range("M56").Value = -6344.5
$ws.Range("H92").Value = 1488.421
$ws.Range("J92").Value = 916.3333
$ws.Range("L92").Value = 2748.9999
$ws.Range("N92").Value = -5244.9999
$ws.Range("H140").Value = 2819.647
$ws.Range("I140").Value = 2819.647
$ws.Range("K140").Value = 8458.940999999999
$ws.Range("M140").Value = -3278.940999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 608.85
$ws.Range("I2").Value = 876.4545000000001
$ws.Range("J2").Value = 281.77777
$ws.Range("K2").Value = 876.4545000000001
$ws.Range("L2").Value = 281.77777
$ws.Range("M2").Value = -763.4545000000001
$ws.Range("N2").Value = -507.77777
$ws.Range("H70").Value = 7091.9116
$ws.Range("J70").Value = 7834.2
$ws.Range("L70").Value = 7834.2
$ws.Range("N70").Value = -8374.200000000001
$ws.Range("H73").Value = 7091.9116
$ws.Range("J73").Value = 7834.2
$ws.Range("L73").Value = 7834.2
$ws.Range("N73").Value = -9706.200000000001
$ws.Range("H126").Value = 3288.1177
$ws.Range("I126").Value = 2100.111
$ws.Range("J126").Value = 4624.625
$ws.Range("K126").Value = 6300.333
$ws.Range("L126").Value = 13873.875
$ws.Range("M126").Value = -3830.333
$ws.Range("N126").Value = -18813.875
$ws.Range("H132").Value = 3345.913
$ws.Range("I132").Value = 3369.8096
$ws.Range("J132").Value = 3095
$ws.Range("K132").Value = 10109.4288
$ws.Range("L132").Value = 9285
$ws.Range("M132").Value = -7579.4288
$ws.Range("N132").Value = -14345

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2279.5356
$ws.Range("I40").Value = 1774.2273
$ws.Range("K40").Value = 1774.2273
$ws.Range("M40").Value = -1638.2273
$ws.Range("H68").Value = 3144.1904
$ws.Range("I68").Value = 2668.6667
$ws.Range("J68").Value = 4333
$ws.Range("K68").Value = 2668.6667
$ws.Range("L68").Value = 4333
$ws.Range("M68").Value = -1919.6667
$ws.Range("N68").Value = -5831
$ws.Range("H71").Value = 3144.1904
$ws.Range("I71").Value = 2668.6667
$ws.Range("J71").Value = 4333
$ws.Range("K71").Value = 13343.3335
$ws.Range("L71").Value = 21665
$ws.Range("M71").Value = -9599.333500000001
$ws.Range("N71").Value = -29153
$ws.Range("H100").Value = 5347
$ws.Range("I100").Value = 4875.2
$ws.Range("J100").Value = 6133.3335
$ws.Range("K100").Value = 4875.2
$ws.Range("L100").Value = 6133.3335
$ws.Range("M100").Value = -4334.2
$ws.Range("N100").Value = -7215.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 120000
$ws.Range("J70").Value = 120000
$ws.Range("L70").Value = 120000
$ws.Range("N70").Value = -120630
$ws.Range("H73").Value = 120000
$ws.Range("J73").Value = 120000
$ws.Range("L73").Value = 120000
$ws.Range("N73").Value = -122184
$ws.Range("H107").Value = 712
$ws.Range("I107").Value = 665.1429000000001
$ws.Range("J107").Value = 805.7143
$ws.Range("K107").Value = 1995.4287
$ws.Range("L107").Value = 2417.1429
$ws.Range("M107").Value = -75.42870000000016
$ws.Range("N107").Value = -6257.1429
$ws.Range("H132").Value = 1893.6818
$ws.Range("I132").Value = 1437.3334
$ws.Range("J132").Value = 3947.25
$ws.Range("K132").Value = 4312.0002
$ws.Range("L132").Value = 11841.75
$ws.Range("M132").Value = -1782.0002
$ws.Range("N132").Value = -16901.75
